$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-01 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-02 Monday", 2)

$d.Content.Find.Execute("348×6=2088", $true, $false, $false, $false, $false, $true, 1, $false, "157×5=785", 2)
$d.Content.Find.Execute("172×7=1204", $true, $false, $false, $false, $false, $true, 1, $false, "630×5=3150", 2)
$d.Content.Find.Execute("822×4=3288", $true, $false, $false, $false, $false, $true, 1, $false, "333×7=2331", 2)
$d.Content.Find.Execute("210×3=630", $true, $false, $false, $false, $false, $true, 1, $false, "326×3=978", 2)
$d.Content.Find.Execute("277×3=831", $true, $false, $false, $false, $false, $true, 1, $false, "392×6=2352", 2)

$d.Content.Find.Execute("224×7=1568", $true, $false, $false, $false, $false, $true, 1, $false, "375×8=3000", 2)
$d.Content.Find.Execute("500×7=3500", $true, $false, $false, $false, $false, $true, 1, $false, "334×6=2004", 2)
$d.Content.Find.Execute("869×4=3476", $true, $false, $false, $false, $false, $true, 1, $false, "174×8=1392", 2)
$d.Content.Find.Execute("544×5=2720", $true, $false, $false, $false, $false, $true, 1, $false, "739×7=5173", 2)
$d.Content.Find.Execute("415×3=1245", $true, $false, $false, $false, $false, $true, 1, $false, "932×9=8388", 2)

$d.Content.Find.Execute("263×2=526", $true, $false, $false, $false, $false, $true, 1, $false, "988×6=5928", 2)
$d.Content.Find.Execute("534×7=3738", $true, $false, $false, $false, $false, $true, 1, $false, "117×6=702", 2)
$d.Content.Find.Execute("890×4=3560", $true, $false, $false, $false, $false, $true, 1, $false, "983×8=7864", 2)
$d.Content.Find.Execute("300×5=1500", $true, $false, $false, $false, $false, $true, 1, $false, "914×9=8226", 2)
$d.Content.Find.Execute("309×7=2163", $true, $false, $false, $false, $false, $true, 1, $false, "544×3=1632", 2)

$d.Content.Find.Execute("548×7=3836", $true, $false, $false, $false, $false, $true, 1, $false, "204×6=1224", 2)
$d.Content.Find.Execute("798×9=7182", $true, $false, $false, $false, $false, $true, 1, $false, "495×3=1485", 2)
$d.Content.Find.Execute("833×6=4998", $true, $false, $false, $false, $false, $true, 1, $false, "102×2=204", 2)
$d.Content.Find.Execute("698×9=6282", $true, $false, $false, $false, $false, $true, 1, $false, "901×9=8109", 2)
$d.Content.Find.Execute("394×8=3152", $true, $false, $false, $false, $false, $true, 1, $false, "815×3=2445", 2)

$d.Content.Find.Execute("651×7=4557", $true, $false, $false, $false, $false, $true, 1, $false, "997×6=5982", 2)
$d.Content.Find.Execute("841×4=3364", $true, $false, $false, $false, $false, $true, 1, $false, "238×3=714", 2)
$d.Content.Find.Execute("324×6=1944", $true, $false, $false, $false, $false, $true, 1, $false, "214×7=1498", 2)
$d.Content.Find.Execute("389×2=778", $true, $false, $false, $false, $false, $true, 1, $false, "874×9=7866", 2)
$d.Content.Find.Execute("518×6=3108", $true, $false, $false, $false, $false, $true, 1, $false, "678×5=3390", 2)
